$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cell, $value)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") "43.272.08"
$ws.Range("E2").Value = "  +2.99%  "

Set-TextValue $ws.Range("D3") "2.315.46"
$ws.Range("E3").Value = "  +2.55%  "

Set-TextValue $ws.Range("D4") "0.999"
$ws.Range("E4").Value = "  -0.02%  "

Set-TextValue $ws.Range("D5") "310.40"
$ws.Range("E5").Value = "  +1.80%  "

Set-TextValue $ws.Range("D6") "101.57"
$ws.Range("E6").Value = "  +5.61%  "

Set-TextValue $ws.Range("D7") "0.536"
$ws.Range("E7").Value = "  +2.21%  "

$ws.Range("E8").Value = "  +0.06%  "

$ws.Range("E9").Value = "  +7.06%  "

Set-TextValue $ws.Range("D10") "36.09"
$ws.Range("E10").Value = "  +3.42%  "

Set-TextValue $ws.Range("D11") "0.0816"
$ws.Range("E11").Value = "  +3.47%  "

Set-TextValue $ws.Range("D13") "7.07"
$ws.Range("E13").Value = "  +3.80%  "

Set-TextValue $ws.Range("D14") "2.672.64"
$ws.Range("E14").Value = "  +2.53%  "

Set-TextValue $ws.Range("D15") "14.99"
$ws.Range("E15").Value = "  +3.99%  "

Set-TextValue $ws.Range("D16") "2.311.87"
$ws.Range("E16").Value = "  +2.45%  "

$ws.Range("E17").Value = "  +2.96%  "

Set-TextValue $ws.Range("D18") "43.105.61"
$ws.Range("E18").Value = "  +2.90%  "

Set-TextValue $ws.Range("D19") "12.57"
$ws.Range("E19").Value = "  +1.94%  "

Set-TextValue $ws.Range("D20") "0.0₃0921"
$ws.Range("E20").Value = "  +2.28%  "

Set-TextValue $ws.Range("D21") "6.14"
$ws.Range("E21").Value = "  +3.21%  "

Set-TextValue $ws.Range("D22") "68.47"
$ws.Range("E22").Value = "  +0.07%  "

Set-TextValue $ws.Range("D23") "241.46"
$ws.Range("E23").Value = "  +1.75%  "

$ws.Range("E24").Value = "  +6.32%  "

Set-TextValue $ws.Range("D25") "2.64"
$ws.Range("E25").Value = "  +3.31%  "

Set-TextValue $ws.Range("D26") "0.999"
$ws.Range("E26").Value = "  -0.19%  "

Set-TextValue $ws.Range("D27") "24.75"
$ws.Range("E27").Value = "  +5.13%  "

Set-TextValue $ws.Range("D28") "37.60"
$ws.Range("E28").Value = "  +3.05%  "

Set-TextValue $ws.Range("D29") "9.67"
$ws.Range("E29").Value = "  +2.68%  "

Set-TextValue $ws.Range("D30") "2.12"
$ws.Range("E30").Value = "  -0.17%  "

Set-TextValue $ws.Range("D31") "167.47"
$ws.Range("E31").Value = "  +4.42%  "

$ws.Range("E32").Value = "  +2.73%  "

Set-TextValue $ws.Range("D33") "0.998"
$ws.Range("E33").Value = "  -0.11%  "

$ws.Range("E34").Value = "  +0.22%  "

Set-TextValue $ws.Range("D35") "17.97"
$ws.Range("E35").Value = "  +5.85%  "

Set-TextValue $ws.Range("D36") "0.0745"
$ws.Range("E36").Value = "  +1.53%  "

$ws.Range("E37").Value = "  +2.71%  "

Set-TextValue $ws.Range("D38") "2.40"
$ws.Range("E38").Value = "  +0.92%  "

Set-TextValue $ws.Range("D39") "1.86"
$ws.Range("E39").Value = "  +3.12%  "

$ws.Range("E40").Value = "  +2.17%  "

$ws.Range("E41").Value = "  +7.99%  "

Set-TextValue $ws.Range("D42") "19.88"
$ws.Range("E42").Value = "  +7.20%  "

$ws.Range("E43").Value = "  +0.97%  "

Set-TextValue $ws.Range("D44") "0.0291"
$ws.Range("E44").Value = "  +3.41%  "

Set-TextValue $ws.Range("D45") "1.975.95"
$ws.Range("E45").Value = "  +0.80%  "

Set-TextValue $ws.Range("D46") "3.01"
$ws.Range("E46").Value = "  +4.28%  "

$ws.Range("E47").Value = "  -1.31%  "

Set-TextValue $ws.Range("D48") "2.99"
$ws.Range("E48").Value = "  +19.03%  "

$ws.Range("E49").Value = "  +4.96%  "

Set-TextValue $ws.Range("D50") "2.541.38"
$ws.Range("E50").Value = "  +2.58%  "

Set-TextValue $ws.Range("D51") "1.54"
$ws.Range("E51").Value = "  +4.65%  "
